$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 6498.3335  # H40
$ws.Cells.Item(40, 10).Value = 6498.3335  # J40
$ws.Cells.Item(40, 12).Value = 6498.3335  # L40
$ws.Cells.Item(40, 14).Value = -6848.3335  # N40
$ws.Cells.Item(55, 8).Value = 90.40000000000001  # H55
$ws.Cells.Item(55, 10).Value = 127.5  # J55
$ws.Cells.Item(55, 12).Value = 127.5  # L55
$ws.Cells.Item(55, 14).Value = -555.5  # N55
$ws.Cells.Item(64, 8).Value = 30311602  # H64
$ws.Cells.Item(64, 9).Value = 111114210  # I64
$ws.Cells.Item(64, 10).Value = 10624.875  # J64
$ws.Cells.Item(64, 11).Value = 111114210  # K64
$ws.Cells.Item(64, 12).Value = 10624.875  # L64
$ws.Cells.Item(64, 13).Value = -111113962  # M64
$ws.Cells.Item(64, 14).Value = -11120.875  # N64
$ws.Cells.Item(67, 8).Value = 30311602  # H67
$ws.Cells.Item(67, 9).Value = 111114210  # I67
$ws.Cells.Item(67, 10).Value = 10624.875  # J67
$ws.Cells.Item(67, 11).Value = 111114210  # K67
$ws.Cells.Item(67, 12).Value = 10624.875  # L67
$ws.Cells.Item(67, 13).Value = -111113352  # M67
$ws.Cells.Item(67, 14).Value = -12340.875  # N67
$ws.Cells.Item(74, 8).Value = 20414834  # H74
$ws.Cells.Item(74, 9).Value = 28574768  # I74
$ws.Cells.Item(74, 10).Value = 15000  # J74
$ws.Cells.Item(74, 11).Value = 28574768  # K74
$ws.Cells.Item(74, 12).Value = 15000  # L74
$ws.Cells.Item(74, 13).Value = -28573832  # M74
$ws.Cells.Item(74, 14).Value = -16872  # N74
$ws.Cells.Item(77, 8).Value = 20414834  # H77
$ws.Cells.Item(77, 9).Value = 28574768  # I77
$ws.Cells.Item(77, 10).Value = 15000  # J77
$ws.Cells.Item(77, 11).Value = 142873840  # K77
$ws.Cells.Item(77, 12).Value = 75000  # L77
$ws.Cells.Item(77, 13).Value = -142869160  # M77
$ws.Cells.Item(77, 14).Value = -84360  # N77
$ws.Cells.Item(88, 8).Value = 3244.1428  # H88
$ws.Cells.Item(88, 9).Value = 2232.4  # I88
$ws.Cells.Item(88, 11).Value = 2232.4  # K88
$ws.Cells.Item(88, 13).Value = -1826.4  # M88
$ws.Cells.Item(91, 8).Value = 3244.1428  # H91
$ws.Cells.Item(91, 9).Value = 2232.4  # I91
$ws.Cells.Item(91, 11).Value = 2232.4  # K91
$ws.Cells.Item(91, 13).Value = -828.4000000000001  # M91
$ws.Cells.Item(100, 8).Value = 2931  # H100
$ws.Cells.Item(100, 9).Value = 1920.2858  # I100
$ws.Cells.Item(100, 11).Value = 1920.2858  # K100
$ws.Cells.Item(100, 13).Value = -1379.2858  # M100
$ws.Cells.Item(112, 8).Value = 3635.0322  # H112
$ws.Cells.Item(112, 10).Value = 3203.1072  # J112
$ws.Cells.Item(112, 12).Value = 9609.321599999999  # L112
$ws.Cells.Item(112, 14).Value = -11825.3216  # N112
$ws.Cells.Item(125, 8).Value = 2398.3333  # H125
$ws.Cells.Item(125, 9).Value = 2966.5  # I125
$ws.Cells.Item(125, 11).Value = 26698.5  # K125
$ws.Cells.Item(125, 13).Value = -24238.5  # M125
$ws.Cells.Item(132, 8).Value = 242508.14  # H132
$ws.Cells.Item(132, 9).Value = 284260.8  # I132
$ws.Cells.Item(132, 10).Value = 5909.6665  # J132
$ws.Cells.Item(132, 11).Value = 852782.3999999999  # K132
$ws.Cells.Item(132, 12).Value = 17728.9995  # L132
$ws.Cells.Item(132, 13).Value = -850252.3999999999  # M132
$ws.Cells.Item(132, 14).Value = -22788.9995  # N132
$ws.Cells.Item(137, 8).Value = 5287.36  # H137
$ws.Cells.Item(137, 9).Value = 5116.5  # I137
$ws.Cells.Item(137, 11).Value = 15349.5  # K137
$ws.Cells.Item(137, 13).Value = -12799.5  # M137
$ws.Cells.Item(138, 8).Value = 4179.656  # H138
$ws.Cells.Item(138, 10).Value = 5363.3823  # J138
$ws.Cells.Item(138, 12).Value = 16090.1469  # L138
$ws.Cells.Item(138, 14).Value = -26370.1469  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1628.4286  # H45
$ws.Cells.Item(45, 10).Value = 3350  # J45
$ws.Cells.Item(45, 12).Value = 3350  # L45
$ws.Cells.Item(45, 14).Value = -4104  # N45
$ws.Cells.Item(56, 8).Value = 4453.846  # H56
$ws.Cells.Item(56, 9).Value = 4453.846  # I56
$ws.Cells.Item(56, 11).Value = 4453.846  # K56
$ws.Cells.Item(56, 13).Value = -3711.846  # M56
$ws.Cells.Item(74, 8).Value = 5489.0835  # H74
$ws.Cells.Item(74, 9).Value = 6353.923  # I74
$ws.Cells.Item(74, 11).Value = 6353.923  # K74
$ws.Cells.Item(74, 13).Value = -5479.923  # M74
$ws.Cells.Item(77, 8).Value = 5489.0835  # H77
$ws.Cells.Item(77, 9).Value = 6353.923  # I77
$ws.Cells.Item(77, 11).Value = 31769.615  # K77
$ws.Cells.Item(77, 13).Value = -27401.615  # M77
$ws.Cells.Item(97, 8).Value = 2132.6667  # H97
$ws.Cells.Item(97, 9).Value = 2199.5  # I97
$ws.Cells.Item(97, 11).Value = 2199.5  # K97
$ws.Cells.Item(97, 13).Value = -1703.5  # M97
$ws.Cells.Item(132, 8).Value = 706648.25  # H132
$ws.Cells.Item(132, 9).Value = 877399.9  # I132
$ws.Cells.Item(132, 10).Value = 72428.07000000001  # J132
$ws.Cells.Item(132, 11).Value = 2632199.7  # K132
$ws.Cells.Item(132, 12).Value = 217284.21  # L132
$ws.Cells.Item(132, 13).Value = -2629669.7  # M132
$ws.Cells.Item(132, 14).Value = -222344.21  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 8387.964  # H99
$ws.Cells.Item(99, 9).Value = 7743.6924  # I99
$ws.Cells.Item(99, 10).Value = 9434.906000000001  # J99
$ws.Cells.Item(99, 11).Value = 7743.6924  # K99
$ws.Cells.Item(99, 12).Value = 9434.906000000001  # L99
$ws.Cells.Item(99, 13).Value = -6245.6924  # M99
$ws.Cells.Item(99, 14).Value = -12430.906  # N99

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6056.8945  # H31
$ws.Cells.Item(31, 9).Value = 1868.1666  # I31
$ws.Cells.Item(31, 10).Value = 13237.571  # J31
$ws.Cells.Item(31, 11).Value = 1868.1666  # K31
$ws.Cells.Item(31, 12).Value = 13237.571  # L31
$ws.Cells.Item(31, 13).Value = -1573.1666  # M31
$ws.Cells.Item(31, 14).Value = -13827.571  # N31
$ws.Cells.Item(34, 8).Value = 6056.8945  # H34
$ws.Cells.Item(34, 9).Value = 1868.1666  # I34
$ws.Cells.Item(34, 10).Value = 13237.571  # J34
$ws.Cells.Item(34, 11).Value = 1868.1666  # K34
$ws.Cells.Item(34, 12).Value = 13237.571  # L34
$ws.Cells.Item(34, 13).Value = -1666.1666  # M34
$ws.Cells.Item(34, 14).Value = -13641.571  # N34
$ws.Cells.Item(58, 8).Value = 8467.471  # H58
$ws.Cells.Item(58, 9).Value = 3615.2  # I58
$ws.Cells.Item(58, 10).Value = 10489.25  # J58
$ws.Cells.Item(58, 11).Value = 3615.2  # K58
$ws.Cells.Item(58, 12).Value = 10489.25  # L58
$ws.Cells.Item(58, 13).Value = -3412.2  # M58
$ws.Cells.Item(58, 14).Value = -10895.25  # N58
$ws.Cells.Item(62, 8).Value = 5611.6313  # H62
$ws.Cells.Item(62, 9).Value = 5274.6  # I62
$ws.Cells.Item(62, 11).Value = 5274.6  # K62
$ws.Cells.Item(62, 13).Value = -4650.6  # M62
$ws.Cells.Item(65, 8).Value = 5611.6313  # H65
$ws.Cells.Item(65, 9).Value = 5274.6  # I65
$ws.Cells.Item(65, 11).Value = 26373  # K65
$ws.Cells.Item(65, 13).Value = -23253  # M65
$ws.Cells.Item(86, 8).Value = 6382.353  # H86
$ws.Cells.Item(86, 9).Value = 5678.1665  # I86
$ws.Cells.Item(86, 10).Value = 7174.5625  # J86
$ws.Cells.Item(86, 11).Value = 5678.1665  # K86
$ws.Cells.Item(86, 12).Value = 7174.5625  # L86
$ws.Cells.Item(86, 13).Value = -4555.1665  # M86
$ws.Cells.Item(86, 14).Value = -9420.5625  # N86
$ws.Cells.Item(89, 8).Value = 6382.353  # H89
$ws.Cells.Item(89, 9).Value = 5678.1665  # I89
$ws.Cells.Item(89, 10).Value = 7174.5625  # J89
$ws.Cells.Item(89, 11).Value = 28390.8325  # K89
$ws.Cells.Item(89, 12).Value = 35872.8125  # L89
$ws.Cells.Item(89, 13).Value = -22774.8325  # M89
$ws.Cells.Item(89, 14).Value = -47104.8125  # N89
$ws.Cells.Item(99, 8).Value = 5295133.5  # H99
$ws.Cells.Item(99, 9).Value = 12349466  # I99
$ws.Cells.Item(99, 10).Value = 4384.1665  # J99
$ws.Cells.Item(99, 11).Value = 12349466  # K99
$ws.Cells.Item(99, 12).Value = 4384.1665  # L99
$ws.Cells.Item(99, 13).Value = -12347968  # M99
$ws.Cells.Item(99, 14).Value = -7380.1665  # N99
$ws.Cells.Item(126, 8).Value = 5295133.5  # H126
$ws.Cells.Item(126, 9).Value = 12349466  # I126
$ws.Cells.Item(126, 10).Value = 4384.1665  # J126
$ws.Cells.Item(126, 11).Value = 37048398  # K126
$ws.Cells.Item(126, 12).Value = 13152.4995  # L126
$ws.Cells.Item(126, 13).Value = -37045928  # M126
$ws.Cells.Item(126, 14).Value = -18092.4995  # N126
$ws.Cells.Item(136, 8).Value = 8467.471  # H136
$ws.Cells.Item(136, 9).Value = 3615.2  # I136
$ws.Cells.Item(136, 10).Value = 10489.25  # J136
$ws.Cells.Item(136, 11).Value = 10845.6  # K136
$ws.Cells.Item(136, 12).Value = 31467.75  # L136
$ws.Cells.Item(136, 13).Value = -8295.599999999999  # M136
$ws.Cells.Item(136, 14).Value = -36567.75  # N136
$ws.Cells.Item(141, 8).Value = 190672.78  # H141
$ws.Cells.Item(141, 10).Value = 190672.78  # J141
$ws.Cells.Item(141, 12).Value = 190672.78  # L141
$ws.Cells.Item(141, 14).Value = -201032.78  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(122, 8).Value = 2839.6  # H122
$ws.Cells.Item(122, 9).Value = 1170.4  # I122
$ws.Cells.Item(122, 11).Value = 10533.6  # K122
$ws.Cells.Item(122, 13).Value = -8083.6  # M122

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1932.1111  # H97
$ws.Cells.Item(97, 9).Value = 1936.625  # I97
$ws.Cells.Item(97, 11).Value = 1936.625  # K97
$ws.Cells.Item(97, 13).Value = -1440.625  # M97
$ws.Cells.Item(113, 8).Value = 5462.5884  # H113
$ws.Cells.Item(113, 9).Value = 4746.3335  # I113
$ws.Cells.Item(113, 10).Value = 5853.273  # J113
$ws.Cells.Item(113, 11).Value = 4746.3335  # K113
$ws.Cells.Item(113, 12).Value = 5853.273  # L113
$ws.Cells.Item(113, 13).Value = -2576.3335  # M113
$ws.Cells.Item(113, 14).Value = -10193.273  # N113
$ws.Cells.Item(132, 8).Value = 6267.0303  # H132
$ws.Cells.Item(132, 9).Value = 6860.92  # I132
$ws.Cells.Item(132, 10).Value = 4411.125  # J132
$ws.Cells.Item(132, 11).Value = 20582.76  # K132
$ws.Cells.Item(132, 12).Value = 13233.375  # L132
$ws.Cells.Item(132, 13).Value = -18052.76  # M132
$ws.Cells.Item(132, 14).Value = -18293.375  # N132
$ws.Cells.Item(135, 8).Value = 111110.75  # H135
$ws.Cells.Item(135, 10).Value = 111110.75  # J135
$ws.Cells.Item(135, 12).Value = 111110.75  # L135
$ws.Cells.Item(135, 14).Value = -121250.75  # N135

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 538.62067  # H22
$ws.Cells.Item(22, 9).Value = 452  # I22
$ws.Cells.Item(22, 10).Value = 680.36365  # J22
$ws.Cells.Item(22, 11).Value = 452  # K22
$ws.Cells.Item(22, 12).Value = 680.36365  # L22
$ws.Cells.Item(22, 13).Value = -157  # M22
$ws.Cells.Item(22, 14).Value = -1270.36365  # N22
$ws.Cells.Item(27, 8).Value = 538.62067  # H27
$ws.Cells.Item(27, 9).Value = 452  # I27
$ws.Cells.Item(27, 10).Value = 680.36365  # J27
$ws.Cells.Item(27, 11).Value = 452  # K27
$ws.Cells.Item(27, 12).Value = 680.36365  # L27
$ws.Cells.Item(27, 13).Value = -345  # M27
$ws.Cells.Item(27, 14).Value = -894.36365  # N27

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(24, 8).Value = 0  # H24
$ws.Cells.Item(24, 10).Value = 0  # J24
$ws.Cells.Item(24, 12).Value = 0  # L24
$ws.Cells.Item(24, 14).ClearContents()  # N24
$ws.Cells.Item(51, 8).Value = 40000  # H51
$ws.Cells.Item(51, 9).Value = 0  # I51
$ws.Cells.Item(51, 11).Value = 0  # K51
$ws.Cells.Item(51, 13).ClearContents()  # M51
$ws.Cells.Item(81, 8).Value = 2429.0667  # H81
$ws.Cells.Item(81, 10).Value = 4187.25  # J81
$ws.Cells.Item(81, 12).Value = 8374.5  # L81
$ws.Cells.Item(81, 14).Value = -10496.5  # N81
$ws.Cells.Item(84, 8).Value = 2429.0667  # H84
$ws.Cells.Item(84, 10).Value = 4187.25  # J84
$ws.Cells.Item(84, 12).Value = 41872.5  # L84
$ws.Cells.Item(84, 14).Value = -52480.5  # N84
$ws.Cells.Item(96, 8).Value = 3558.4  # H96
$ws.Cells.Item(96, 9).Value = 3558.4  # I96
$ws.Cells.Item(96, 10).Value = 0  # J96
$ws.Cells.Item(96, 11).Value = 3558.4  # K96
$ws.Cells.Item(96, 12).Value = 0  # L96
$ws.Cells.Item(96, 13).Value = -2185.4  # M96
$ws.Cells.Item(96, 14).ClearContents()  # N96
$ws.Cells.Item(136, 8).Value = 21757586  # H136
$ws.Cells.Item(136, 9).Value = 33352336  # I136
$ws.Cells.Item(136, 11).Value = 100057008  # K136
$ws.Cells.Item(136, 13).Value = -100054458  # M136

